$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 5, Target cluster = "Resolving-Mac") entirely
$ws.Rows.Item(5).Delete()

# Update row 2 (Target cluster = ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1127976666666667
$ws.Range("H2").Value = 0.338393
$ws.Range("M2").Value = 30.58864766666666
$ws.Range("N2").Value = 91.76594299999999
$ws.Range("O2").Value = 0.3925391465174898
$ws.Range("P2").Value = 0.3925391465174898
$ws.Range("Q2").Value = 3.450328083288777
$ws.Range("R2").Value = 31.052952749599
$ws.Range("S2").Value = 0.3925391465174898
$ws.Range("T2").Value = 0.3925391465174898

# Update row 3 (Target cluster = FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1127976666666667
$ws.Range("H3").Value = 0.338393
$ws.Range("O3").Value = 0.291183949679193
$ws.Range("P3").Value = 0.291183949679193
$ws.Range("Q3").Value = 2.559439403418334
$ws.Range("R3").Value = 23.034954630765
$ws.Range("S3").Value = 0.291183949679193
$ws.Range("T3").Value = 0.291183949679193

# Update row 4 (Target cluster = MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1127976666666667
$ws.Range("H4").Value = 0.338393
$ws.Range("M4").Value = 24.64590566666666
$ws.Range("N4").Value = 73.93771699999999
$ws.Range("O4").Value = 0.3162769038033173
$ws.Range("P4").Value = 0.3162769038033172
$ws.Range("Q4").Value = 2.780000652086778
$ws.Range("R4").Value = 25.020005868781
$ws.Range("S4").Value = 0.3162769038033173
$ws.Range("T4").Value = 0.3162769038033172
